$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "24×80=1920"; New = "41×21=861" },
    @{ Old = "32×41=1312"; New = "18×25=450" },
    @{ Old = "72×38=2736"; New = "34×52=1768" },
    @{ Old = "60×24=1440"; New = "95×67=6365" },
    @{ Old = "65×65=4225"; New = "27×89=2403" },
    @{ Old = "32×89=2848"; New = "28×66=1848" },
    @{ Old = "93×63=5859"; New = "79×34=2686" },
    @{ Old = "47×52=2444"; New = "64×88=5632" },
    @{ Old = "45×47=2115"; New = "51×98=4998" },
    @{ Old = "78×41=3198"; New = "45×73=3285" },
    @{ Old = "95×42=3990"; New = "99×55=5445" },
    @{ Old = "52×40=2080"; New = "48×47=2256" },
    @{ Old = "17×52=884";  New = "39×65=2535" },
    @{ Old = "24×85=2040"; New = "13×69=897" },
    @{ Old = "86×14=1204"; New = "95×70=6650" },
    @{ Old = "39×67=2613"; New = "61×43=2623" },
    @{ Old = "23×56=1288"; New = "88×50=4400" },
    @{ Old = "48×40=1920"; New = "74×34=2516" },
    @{ Old = "34×23=782";  New = "83×79=6557" },
    @{ Old = "28×56=1568"; New = "73×93=6789" },
    @{ Old = "26×55=1430"; New = "98×65=6370" },
    @{ Old = "87×29=2523"; New = "50×64=3200" },
    @{ Old = "59×78=4602"; New = "76×94=7144" },
    @{ Old = "53×32=1696"; New = "56×82=4592" },
    @{ Old = "46×23=1058"; New = "23×96=2208" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
